$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the "EFS"/"OS" column headers in B1/C1 (route through a temp value so
# the two assignments don't clobber each other via shared-string aliasing).
$ws.Range("B1").Value = "__TMP_HEADER_SWAP__"
$ws.Range("C1").Value = "Significant CpG probes in EFS"
$ws.Range("B1").Value = "Significant CpG probes in OS"

# Copy the "risk group" label formatting (bold, border, centered) from A2
# down into A3 before we repurpose A2's text.
$ws.Range("A2").Copy()
$ws.Range("A3").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 2 becomes "Previous Risk Group" with its own EWAS counts.
$ws.Range("A2").Value = "Previous Risk Group"
$ws.Range("B2").Value = 609
$ws.Range("C2").Value = 256
$ws.Range("D2").Value = 77

# New row 3 holds the "Updated Risk Group" counts (previously on row 2).
$ws.Range("A3").Value = "Updated Risk Group"
$ws.Range("B3").Value = 189
$ws.Range("C3").Value = 112
$ws.Range("D3").Value = 17
